# Remove the "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph and the
# "© 2020 . Contact: ..." footer paragraph that follow the Requisitos section,
# while leaving the surrounding blank paragraphs untouched.

$d = $word.ActiveDocument

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = [char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$p1 = $null
$p2 = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq $target1) {
        $p1 = $p
    } elseif ($text -eq $target2) {
        $p2 = $p
    }
}

$start = $p1.Range.Start
$end = $p2.Range.End
$r = $d.Range($start, $end)
$r.Delete()
